$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writing a value like "304.80" or "2.22%" via .Value alone lets Excel's
# usual type inference kick in (it gets parsed as a number / percentage and
# picks up a new number-format style). The source data keeps these columns
# as plain text, so force text entry with a leading apostrophe and then
# reset the cell style back to "Normal" so no extra style gets attached.
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# New Price (D), Volume/1h (E) and Hora (G) values for rows 2-51, taken from
# the refreshed symbol list. $null means that column did not change for that
# row (e.g. rows whose Price/Volume are still the "--" placeholder).
$updates = @(
    @{ Row=2; D='304.80'; E='2.22%'; G='5' },
    @{ Row=3; D='31.77'; E='0.17%'; G='5' },
    @{ Row=4; D='5.170'; E='1.49%'; G='5' },
    @{ Row=5; D=$null; E='-0.31%'; G='5' },
    @{ Row=6; D='2.343'; E='37.76%'; G='5' },
    @{ Row=7; D='8.018'; E='3.05%'; G='5' },
    @{ Row=8; D='3.876'; E='2.08%'; G='5' },
    @{ Row=9; D='0.9163'; E='-1.13%'; G='5' },
    @{ Row=10; D=$null; E='1.66%'; G='5' },
    @{ Row=11; D='0.07677'; E='5.18%'; G='5' },
    @{ Row=12; D='0.08213'; E='3.29%'; G='5' },
    @{ Row=13; D='0.03040'; E='-0.57%'; G='5' },
    @{ Row=14; D='0.09954'; E='0.66%'; G='5' },
    @{ Row=15; D='0.001504'; E='0.02%'; G='5' },
    @{ Row=16; D='0.006160'; E='-4.99%'; G='5' },
    @{ Row=17; D='3.501'; E='1.46%'; G='5' },
    @{ Row=18; D=$null; E=$null; G='5' },
    @{ Row=19; D=$null; E='-0.85%'; G='5' },
    @{ Row=20; D=$null; E='-0.55%'; G='5' },
    @{ Row=21; D='4.651'; E='1.87%'; G='5' },
    @{ Row=22; D='0.04609'; E='-0.95%'; G='5' },
    @{ Row=23; D='0.1563'; E='0.85%'; G='5' },
    @{ Row=24; D=$null; E='3.76%'; G='5' },
    @{ Row=25; D='0.004539'; E='2.73%'; G='5' },
    @{ Row=26; D='0.0001298'; E='-7.37%'; G='5' },
    @{ Row=27; D='0.0002738'; E='48.53%'; G='5' },
    @{ Row=28; D=$null; E=$null; G='5' },
    @{ Row=29; D=$null; E=$null; G='5' },
    @{ Row=30; D=$null; E=$null; G='5' },
    @{ Row=31; D=$null; E=$null; G='5' },
    @{ Row=32; D=$null; E=$null; G='5' },
    @{ Row=33; D=$null; E=$null; G='5' },
    @{ Row=34; D=$null; E=$null; G='5' },
    @{ Row=35; D=$null; E=$null; G='5' },
    @{ Row=36; D=$null; E=$null; G='5' },
    @{ Row=37; D=$null; E=$null; G='5' },
    @{ Row=38; D=$null; E=$null; G='5' },
    @{ Row=39; D='0.01756'; E='4.89%'; G='5' },
    @{ Row=40; D='0.04554'; E='0.14%'; G='5' },
    @{ Row=41; D='0.007379'; E='4.74%'; G='5' },
    @{ Row=42; D='0.1365'; E='2.98%'; G='5' },
    @{ Row=43; D='0.002127'; E='3.14%'; G='5' },
    @{ Row=44; D='0.01091'; E='-15.46%'; G='5' },
    @{ Row=45; D='0.00006495'; E='7.97%'; G='5' },
    @{ Row=46; D=$null; E='-57.48%'; G='5' },
    @{ Row=47; D=$null; E='-23.88%'; G='5' },
    @{ Row=48; D=$null; E=$null; G='5' },
    @{ Row=49; D=$null; E=$null; G='5' },
    @{ Row=50; D=$null; E=$null; G='5' },
    @{ Row=51; D=$null; E=$null; G='5' }
)

foreach ($u in $updates) {
    $rowNum = $u.Row
    if ($null -ne $u.D) { Set-TextValue $ws.Range("D$rowNum") $u.D }
    if ($null -ne $u.E) { Set-TextValue $ws.Range("E$rowNum") $u.E }
    if ($null -ne $u.G) { Set-TextValue $ws.Range("G$rowNum") $u.G }
}
